$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (12) - this shifts SmediumHealing..RTLastSaved right by one
$ws.Columns("L").Insert()

# Set header for the new column
$ws.Range("L1").Value = "Intelligence"

# Set the data row value for the new column
$ws.Range("L2").Value = 0

# Recompute the best-fit width for the newly inserted column
$ws.Columns("L").AutoFit()

# Update selection to match the recorded state
$ws.Range("L3").Select()
